# Updated cryptos list (prices in column D, 1h volume % in column E).
# Source workbook stores every D/E value as literal text (not numbers),
# including values that happen to look numeric (e.g. "608.51") and
# percentage strings padded with spaces (e.g. "  -1.04%  "). Whenever the
# new text would otherwise be auto-parsed by Excel as a number, the cell is
# first switched to the "Text" number format so the literal string is kept.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cell reference -> new text value
$updates = [ordered]@{
    'D2'  = '66.335.40'
    'E2'  = '  -1.04%  '
    'D3'  = '3.537.94'
    'E4'  = '  -0.02%  '
    'D5'  = '608.51'
    'E5'  = '  -0.19%  '
    'D6'  = '144.38'
    'E6'  = '  -2.39%  '
    'E7'  = '  +0.50%  '
    'E8'  = '  +0.11%  '
    'E9'  = '  +0.40%  '
    'D10' = '8.16'
    'E10' = '  +0.82%  '
    'E11' = '  -4.29%  '
    'E12' = '  -2.83%  '
    'D13' = '4.137.83'
    'E13' = '  +0.59%  '
    'E14' = '  -4.62%  '
    'D15' = '30.29'
    'E15' = '  -5.02%  '
    'D16' = '3.537.10'
    'E16' = '  +0.39%  '
    'D17' = '66.406.00'
    'E17' = '  -0.99%  '
    'E18' = '  -0.78%  '
    'D19' = '10.98'
    'E19' = '  +1.53%  '
    'D20' = '6.22'
    'E20' = '  -3.79%  '
    'D21' = '14.93'
    'D22' = '426.34'
    'E22' = '  -2.63%  '
    'E23' = '  -1.29%  '
    'D24' = '78.92'
    'E24' = '  -0.97%  '
    'D25' = '3.678.71'
    'E25' = '  +0.63%  '
    'E26' = '  +0.04%  '
    'E27' = '  -0.53%  '
    'D28' = '8.16'
    'E28' = '  -1.52%  '
    'D29' = '9.27'
    'E29' = '  -5.29%  '
    'E30' = '  -1.86%  '
    'E31' = '  +0.01%  '
    'E32' = '  -6.96%  '
    'E33' = '  -4.60%  '
    'D34' = '25.30'
    'E34' = '  -1.22%  '
    'D35' = '3.528.10'
    'E35' = '  +0.49%  '
    'E36' = '  -0.03%  '
    'E37' = '  -3.34%  '
    'E38' = '  -3.15%  '
    'E39' = '  -5.69%  '
    'E40' = '  -0.06%  '
    'D41' = '173.49'
    'E41' = '  -1.39%  '
    'E42' = '  -4.54%  '
    'D43' = '5.26'
    'E43' = '  -2.73%  '
    'D44' = '0.893'
    'E44' = '  -0.20%  '
    'E45' = '  -7.32%  '
    'D46' = '45.70'
    'E46' = '  -1.10%  '
    'E47' = '  -3.01%  '
    'D48' = '26.09'
    'E48' = '  -7.37%  '
    'D49' = '2.40'
    'E49' = '  -2.40%  '
    'D50' = '7.13'
    'E50' = '  -4.68%  '
    'E51' = '  -5.31%  '
}

foreach ($ref in $updates.Keys) {
    $newValue = $updates[$ref]
    $cell = $ws.Range($ref)

    # Numeric-looking text (e.g. "608.51", "8.16") would otherwise be
    # auto-converted to a Number by Excel; force the Text format first so
    # it is stored as a plain string, matching the source data.
    if ($newValue -match '^[+-]?\d+(\.\d+)?$') {
        $cell.NumberFormat = '@'
    }

    $cell.Value = $newValue
}
